$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "dog": append a new test record as row 22
# ---------------------------------------------------------------------------
$dog = $wb.Worksheets.Item("dog")

$dog.Cells.Item(22, 1).Value = 45811
$dog.Cells.Item(22, 1).NumberFormat = "m/d/yy"

$dog.Cells.Item(22, 2).Value = "PRESENCE"

$dog.Cells.Item(22, 3).Value = 0.30555555555555558
$dog.Cells.Item(22, 3).NumberFormat = "h:mm"

$dog.Cells.Item(22, 4).Value = 0.48958333333333331
$dog.Cells.Item(22, 4).NumberFormat = "h:mm"

$dog.Cells.Item(22, 5).Value = 13
$dog.Cells.Item(22, 6).Value = 9
$dog.Cells.Item(22, 7).Value = "Sunny, mild"
$dog.Cells.Item(22, 8).Value = $true
$dog.Cells.Item(22, 9).Value = "2 minutes 30 seconds"
$dog.Cells.Item(22, 10).Value = 150
$dog.Cells.Item(22, 11).Value = "Primary sweeps"
$dog.Cells.Item(22, 12).Value = "Worked downhill. Tricky to see once found, good search effort."

# Move the selection on the "dog" sheet the way the author left it.
$dog.Range("J25").Select()

# ---------------------------------------------------------------------------
# Sheet "human": rename the searcher and append four new records
# ---------------------------------------------------------------------------
$human = $wb.Worksheets.Item("human")

$human.Cells.Item(2, 1).Value = "Marc"
$human.Cells.Item(3, 1).Value = "Marc"

# Row 4 - Maisie, first attempt (not found)
$human.Cells.Item(4, 1).Value = "Maisie"
$human.Cells.Item(4, 2).Value = 45811
$human.Cells.Item(4, 2).NumberFormat = "m/d/yy"
$human.Cells.Item(4, 3).Value = 0.61111111111111116
$human.Cells.Item(4, 3).NumberFormat = "h:mm"
$human.Cells.Item(4, 4).Value = "Sunny, cool"
$human.Cells.Item(4, 5).Value = $false
$human.Cells.Item(4, 6).Value = "1 hour"
$human.Cells.Item(4, 7).Value = 3600
$human.Cells.Item(4, 8).Value = "Worked downhill. Found rain moth casings. Found it tiring."

# Row 5 - Maisie, second attempt (found)
$human.Cells.Item(5, 1).Value = "Maisie"
$human.Cells.Item(5, 2).Value = 45811
$human.Cells.Item(5, 2).NumberFormat = "m/d/yy"
$human.Cells.Item(5, 3).Value = 0.67013888888888884
$human.Cells.Item(5, 3).NumberFormat = "h:mm"
$human.Cells.Item(5, 4).Value = "Sunny, cool"
$human.Cells.Item(5, 5).Value = $true
$human.Cells.Item(5, 6).Value = "28 minutes 50 seconds"
$human.Cells.Item(5, 7).Value = 1730
$human.Cells.Item(5, 8).Value = "Worked uphill. Found GPS! Find was 3/4 of the way through, worked faster than the first attempt."

# Row 6 - Micha, spiral search (lucky fast find)
$human.Cells.Item(6, 1).Value = "Micha"
$human.Cells.Item(6, 2).Value = 45812
$human.Cells.Item(6, 2).NumberFormat = "m/d/yy"
$human.Cells.Item(6, 3).Value = 0.4513888888888889
$human.Cells.Item(6, 3).NumberFormat = "h:mm"
$human.Cells.Item(6, 4).Value = "Sunny, cool"
$human.Cells.Item(6, 5).Value = $true
$human.Cells.Item(6, 6).Value = "2 minutes 21 seconds"
$human.Cells.Item(6, 7).Value = 141
$human.Cells.Item(6, 8).Value = "Spiral search. Very fast find- entirely luck based, GPS was <2m from centre."

# Row 7 - Micha, spiral search (more realistic)
$human.Cells.Item(7, 1).Value = "Micha"
$human.Cells.Item(7, 2).Value = 45812
$human.Cells.Item(7, 2).NumberFormat = "m/d/yy"
$human.Cells.Item(7, 3).Value = 0.45833333333333331
$human.Cells.Item(7, 3).NumberFormat = "h:mm"
$human.Cells.Item(7, 4).Value = "Sunny, cool"
$human.Cells.Item(7, 5).Value = $true
$human.Cells.Item(7, 6).Value = "25 minutes 27 seconds"
$human.Cells.Item(7, 7).Value = 1527
$human.Cells.Item(7, 8).Value = "Spiral search. More realistic time and search, GPS <10m from centre."

# Column A narrows now that "Marc Layton" is gone, only short first names remain.
$human.Columns.Item(1).ColumnWidth = 7.7

# "human" becomes the active sheet/tab, with this cell selected.
$human.Range("G12").Select()
$human.Activate()
